$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsSelfReg   = $wb.Worksheets.Item("selfRegistration")
$wsDupEmail  = $wb.Worksheets.Item("selfRegWithDuplicateEmail")

# --- selfRegistration: change available language from DUTCH to English (US) ---
[void]$wsSelfReg.Activate()
$wsSelfReg.Range("D2").Value = "English (US)"
[void]$wsSelfReg.Range("E2").Select()
$wsSelfReg.Columns.Item(4).AutoFit()
$wsSelfReg.PageSetup.Orientation = 1

# --- TestCases: duplicateEmail sign up test now ignored (Runmode N) ---
$wsTestCases.Range("C3").Value = "N"

# TestCases becomes the active/selected tab (was selfRegWithDuplicateEmail before)
[void]$wsTestCases.Activate()
[void]$wsTestCases.Range("C4").Select()

$wb.Save()
